$wb = $excel.ActiveWorkbook

# =========================================================================
# Sheet: Parameters -- selection change B2:B4 -> B2
# (select on the Parameters sheet, then reactivate Example so tabSelected stays put)
# =========================================================================
$wsParameters = $wb.Worksheets.Item("Parameters")
$wsExample = $wb.Worksheets.Item("Example")
$wsParameters.Range("B2").Select()

# =========================================================================
# Sheet: Example -- new wide cohort table (columns G..AY, rows 1..4)
# =========================================================================
$wsExample.Activate()
$ws4 = $wsExample

# Row 1: drop stale per-cell styles on cells that should end up with no explicit style
$ws4.Range("A1").ClearFormats()
$ws4.Range("E1:J1").ClearFormats()

# Row 1: new header cells, written left to right (G1..AX1) so new shared-string
# indices are allocated in the same order as the target workbook
$ws4.Cells.Item(1, 7).Value = "is_in_birth12"
$ws4.Cells.Item(1, 8).Value = "cohort_entry_date_birth12"
$ws4.Cells.Item(1, 9).Value = "cohort_exit_date_birth12"
$ws4.Cells.Item(1, 10).Value = "is_censored_in_birth12"
$ws4.Cells.Item(1, 11).Value = "is_in_birth15"
$ws4.Cells.Item(1, 12).Value = "cohort_entry_date_birth15"
$ws4.Cells.Item(1, 13).Value = "cohort_exit_date_birth15"
$ws4.Cells.Item(1, 14).Value = "is_censored_in_birth15"
$ws4.Cells.Item(1, 15).Value = "is_in_birth24"
$ws4.Cells.Item(1, 16).Value = "cohort_entry_date_birth24"
$ws4.Cells.Item(1, 17).Value = "cohort_exit_date_birth24"
$ws4.Cells.Item(1, 18).Value = "is_censored_in_birth24"
$ws4.Cells.Item(1, 19).Value = "is_in_adolescence"
$ws4.Cells.Item(1, 20).Value = "cohort_entry_date_adolescence"
$ws4.Cells.Item(1, 21).Value = "cohort_exit_date_adolescence"
$ws4.Cells.Item(1, 22).Value = "is_censored_in_adolescence"
$ws4.Cells.Item(1, 23).Value = "is_in_covid_vacc"
$ws4.Cells.Item(1, 24).Value = "cohort_entry_date_covid_vacc"
$ws4.Cells.Item(1, 25).Value = "cohort_exit_date_covid_vacc"
$ws4.Cells.Item(1, 26).Value = "is_censored_in_covid_vacc"
$ws4.Cells.Item(1, 27).Value = "is_in_seasonal2018"
$ws4.Cells.Item(1, 28).Value = "cohort_entry_date_seasonal2018"
$ws4.Cells.Item(1, 29).Value = "cohort_exit_date_seasonal2018"
$ws4.Cells.Item(1, 30).Value = "is_censored_in_seasonal2018"
$ws4.Cells.Item(1, 31).Value = "is_in_seasonal2019"
$ws4.Cells.Item(1, 32).Value = "cohort_entry_date_seasonal2019"
$ws4.Cells.Item(1, 33).Value = "cohort_exit_date_seasonal2019"
$ws4.Cells.Item(1, 34).Value = "is_censored_in_seasonal2019"
$ws4.Cells.Item(1, 35).Value = "is_in_seasonal2020"
$ws4.Cells.Item(1, 36).Value = "cohort_entry_date_seasonal2020"
$ws4.Cells.Item(1, 37).Value = "cohort_exit_date_seasonal2020"
$ws4.Cells.Item(1, 38).Value = "is_censored_in_seasonal2020"
$ws4.Cells.Item(1, 39).Value = "is_in_seasonal2021"
$ws4.Cells.Item(1, 40).Value = "cohort_entry_date_seasonal2021"
$ws4.Cells.Item(1, 41).Value = "cohort_exit_date_seasonal2021"
$ws4.Cells.Item(1, 42).Value = "is_censored_in_seasonal2021"
$ws4.Cells.Item(1, 43).Value = "is_in_seasonal2022"
$ws4.Cells.Item(1, 44).Value = "cohort_entry_date_seasonal2022"
$ws4.Cells.Item(1, 45).Value = "cohort_exit_date_seasonal2022"
$ws4.Cells.Item(1, 46).Value = "is_censored_in_seasonal2022"
$ws4.Cells.Item(1, 47).Value = "is_in_seasonal2023"
$ws4.Cells.Item(1, 48).Value = "cohort_entry_date_seasonal2023"
$ws4.Cells.Item(1, 49).Value = "cohort_exit_date_seasonal2023"
$ws4.Cells.Item(1, 50).Value = "is_censored_in_seasonal2023"

# Data rows
# Row 2 (P001)
$ws4.Cells.Item(2, 1).Value = "P001"
$ws4.Cells.Item(2, 2).Value = 20200101
$ws4.Cells.Item(2, 3).Value = 20200101
$ws4.Cells.Item(2, 4).Value = 20231231
$ws4.Cells.Item(2, 5).Value = 20200101
$ws4.Cells.Item(2, 6).Value = 2020
$ws4.Cells.Item(2, 7).Value = 1
$ws4.Cells.Item(2, 8).Value = 20200101
$ws4.Cells.Item(2, 9).Value = 20201231
$ws4.Cells.Item(2, 10).Value = 0
$ws4.Cells.Item(2, 11).Value = 1
$ws4.Cells.Item(2, 12).Value = 20200101
$ws4.Cells.Item(2, 13).Value = 20210331
$ws4.Cells.Item(2, 14).Value = 0
$ws4.Cells.Item(2, 15).Value = 1
$ws4.Cells.Item(2, 16).Value = 20200101
$ws4.Cells.Item(2, 17).Value = 20211231
$ws4.Cells.Item(2, 18).Value = 0
$ws4.Cells.Item(2, 19).Value = 0
$ws4.Cells.Item(2, 23).Value = 1
$ws4.Cells.Item(2, 24).Value = 20201201
$ws4.Cells.Item(2, 25).Value = 20231231
$ws4.Cells.Item(2, 26).Value = 0
$ws4.Cells.Item(2, 27).Value = 0

# Row 3 (P002)
$ws4.Cells.Item(3, 1).Value = "P002"
$ws4.Cells.Item(3, 2).Value = 20100101
$ws4.Cells.Item(3, 3).Value = 20180101
$ws4.Cells.Item(3, 4).Value = 20231231
$ws4.Cells.Item(3, 5).Value = 20100101
$ws4.Cells.Item(3, 6).Value = 2010
$ws4.Cells.Item(3, 7).Value = 0
$ws4.Cells.Item(3, 11).Value = 0
$ws4.Cells.Item(3, 15).Value = 0
$ws4.Cells.Item(3, 19).Value = 1
$ws4.Cells.Item(3, 20).Value = 20190101
$ws4.Cells.Item(3, 21).Value = 20231231
$ws4.Cells.Item(3, 22).Value = 1
$ws4.Cells.Item(3, 23).Value = 1
$ws4.Cells.Item(3, 24).Value = 20201201
$ws4.Cells.Item(3, 25).Value = 20231231
$ws4.Cells.Item(3, 26).Value = 0
$ws4.Cells.Item(3, 27).Value = 0

# Row 4 (P003)
$ws4.Cells.Item(4, 1).Value = "P003"
$ws4.Cells.Item(4, 2).Value = 19400101
$ws4.Cells.Item(4, 3).Value = 20180101
$ws4.Cells.Item(4, 4).Value = 20231231
$ws4.Cells.Item(4, 5).Value = 19400101
$ws4.Cells.Item(4, 6).Value = 1940
$ws4.Cells.Item(4, 7).Value = 0
$ws4.Cells.Item(4, 11).Value = 0
$ws4.Cells.Item(4, 15).Value = 0
$ws4.Cells.Item(4, 19).Value = 0
$ws4.Cells.Item(4, 23).Value = 1
$ws4.Cells.Item(4, 24).Value = 20201201
$ws4.Cells.Item(4, 25).Value = 20231231
$ws4.Cells.Item(4, 26).Value = 0
$ws4.Cells.Item(4, 27).Value = 1
$ws4.Cells.Item(4, 28).Value = 20180901
$ws4.Cells.Item(4, 29).Value = 20190430
$ws4.Cells.Item(4, 30).Value = 0
$ws4.Cells.Item(4, 31).Value = 1
$ws4.Cells.Item(4, 32).Value = 20190901
$ws4.Cells.Item(4, 33).Value = 20200430
$ws4.Cells.Item(4, 34).Value = 0
$ws4.Cells.Item(4, 35).Value = 1
$ws4.Cells.Item(4, 36).Value = 20200901
$ws4.Cells.Item(4, 37).Value = 20210430
$ws4.Cells.Item(4, 38).Value = 0
$ws4.Cells.Item(4, 39).Value = 1
$ws4.Cells.Item(4, 40).Value = 20210901
$ws4.Cells.Item(4, 41).Value = 20220430
$ws4.Cells.Item(4, 42).Value = 0
$ws4.Cells.Item(4, 43).Value = 1
$ws4.Cells.Item(4, 44).Value = 20220901
$ws4.Cells.Item(4, 45).Value = 20230430
$ws4.Cells.Item(4, 46).Value = 0
$ws4.Cells.Item(4, 47).Value = 0
$ws4.Cells.Item(4, 48).Value = 1
$ws4.Cells.Item(4, 49).Value = 20230901
$ws4.Cells.Item(4, 50).Value = 20231231
$ws4.Cells.Item(4, 51).Value = 1

# Final selection on Example mirrors the authored workbook (full data block selected)
$ws4.Range("A1:AY4").Select()
